$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet and update header text to reflect new "through" date
$ws.Name = "Through 2022-05-16"
$ws.Range("A6").Value = "May (through 05-16)"

# Update May row (row 6) values
$ws.Range("B6").Value = 13
$ws.Range("C6").Value = 24
$ws.Range("D6").Value = 34
$ws.Range("E6").Value = 24
$ws.Range("F6").Value = 23
$ws.Range("G6").Value = 32
$ws.Range("H6").Value = 59
$ws.Range("I6").Value = 59

# Update Total row (row 7) values
$ws.Range("B7").Value = 102
$ws.Range("C7").Value = 186
$ws.Range("D7").Value = 287
$ws.Range("E7").Value = 270
$ws.Range("F7").Value = 178
$ws.Range("G7").Value = 294
$ws.Range("H7").Value = 582
$ws.Range("I7").Value = 611
